$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '243.27'
$ws.Range('D2').NumberFormat = "General"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '23.04'
$ws.Range('D3').NumberFormat = "General"

$ws.Range('B4').Value = 'HuobiToken'

$ws.Range('C4').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.392'
$ws.Range('D4').NumberFormat = "General"

$ws.Range('E4').Value = '3HuobiTokenHT'

$ws.Range('B5').Value = 'Cronos'

$ws.Range('C5').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05913'
$ws.Range('D5').NumberFormat = "General"

$ws.Range('E5').Value = '4CronosCRO'

$ws.Range('B6').Value = 'GateToken'

$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '3.456'
$ws.Range('D6').NumberFormat = "General"

$ws.Range('E6').Value = '5GateTokenGT'

$ws.Range('B7').Value = 'KuCoinToken'

$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.559'
$ws.Range('D7').NumberFormat = "General"

$ws.Range('E7').Value = '6KuCoinTokenKCS'

$ws.Range('B8').Value = 'MXToken'

$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8108'
$ws.Range('D8').NumberFormat = "General"

$ws.Range('E8').Value = '7MXTokenMX'

$ws.Range('B9').Value = 'FTXToken'

$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9205'
$ws.Range('D9').NumberFormat = "General"

$ws.Range('E9').Value = '8FTXTokenFTT'

$ws.Range('B10').Value = 'WazirX'

$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1412'
$ws.Range('D10').NumberFormat = "General"

$ws.Range('E10').Value = '9WazirXWRX'

$ws.Range('B11').Value = 'MandalaExchangeToken'

$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07425'
$ws.Range('D11').NumberFormat = "General"

$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'

$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.03266'
$ws.Range('D12').NumberFormat = "General"

$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range('B13').Value = 'BitrueCoin'

$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03071'
$ws.Range('D13').NumberFormat = "General"

$ws.Range('E13').Value = '12BitrueCoinBTR'

$ws.Range('B14').Value = 'BitMartToken'

$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.09340'
$ws.Range('D14').NumberFormat = "General"

$ws.Range('E14').Value = '13BitMartTokenBMX'

$ws.Range('B15').Value = 'MCDex'

$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.850'
$ws.Range('D15').NumberFormat = "General"

$ws.Range('E15').Value = '14MCDexMCB'

$ws.Range('B16').Value = 'BitForexToken'

$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.001570'
$ws.Range('D16').NumberFormat = "General"

$ws.Range('E16').Value = '15BitForexTokenBF'

$ws.Range('B17').Value = 'CoinExToken'

$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.04667'
$ws.Range('D17').NumberFormat = "General"

$ws.Range('E17').Value = '16CoinExTokenCET'

$ws.Range('B18').Value = 'One'

$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0005941'
$ws.Range('D18').NumberFormat = "General"

$ws.Range('E18').Value = '17OneONE'

$ws.Range('B19').Value = 'TigerCash'

$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.005919'
$ws.Range('D19').NumberFormat = "General"

$ws.Range('E19').Value = '18TigerCashTCH'

$ws.Range('B20').Value = 'HotbitToken'

$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.004938'
$ws.Range('D20').NumberFormat = "General"

$ws.Range('E20').Value = '19HotbitTokenHTB'

$ws.Range('B21').Value = 'BitKan'

$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0009834'
$ws.Range('D21').NumberFormat = "General"

$ws.Range('E21').Value = '20BitKanKAN'

$ws.Range('B22').Value = 'NitroEx'

$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.00009505'
$ws.Range('D22').NumberFormat = "General"

$ws.Range('E22').Value = '21NitroExNTX'

$ws.Range('B23').Value = 'LEO'

$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.616'
$ws.Range('D23').NumberFormat = "General"

$ws.Range('E23').Value = '22LEOLEO'

$ws.Range('B24').Value = 'BTSEToken'

$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.148'
$ws.Range('D24').NumberFormat = "General"

$ws.Range('E24').Value = '23BTSETokenBTSE'

$ws.Range('B25').Value = 'BitpandaEcosystemToken'

$ws.Range('C25').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.3240'
$ws.Range('D25').NumberFormat = "General"

$ws.Range('E25').Value = '24BitpandaEcosystemTokenBEST'

$ws.Range('B26').Value = 'ProBitToken'

$ws.Range('C26').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1328'
$ws.Range('D26').NumberFormat = "General"

$ws.Range('E26').Value = '25ProBitTokenPROB'

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0002901'
$ws.Range('D27').NumberFormat = "General"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03972'
$ws.Range('D40').NumberFormat = "General"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.006183'
$ws.Range('D41').NumberFormat = "General"

$ws.Range('E41').Value = '40KickTokenKICK'

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1076'
$ws.Range('D42').NumberFormat = "General"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002551'
$ws.Range('D43').NumberFormat = "General"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.008109'
$ws.Range('D44').NumberFormat = "General"

$ws.Range('E44').Value = '43LocalTradersLCTBestin24h'

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005176'
$ws.Range('D45').NumberFormat = "General"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.8402'
$ws.Range('D47').NumberFormat = "General"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.002283'
$ws.Range('D48').NumberFormat = "General"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002100'
$ws.Range('D49').NumberFormat = "General"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0002000'
$ws.Range('D50').NumberFormat = "General"
